$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Output "NOT FOUND: $find"
    }
}

# 1. "launched soon after 2009" -> "launched around 2009"
Replace-Text "was launched soon after 2009" "was launched around 2009"

# 2. quality/reusability wording in paragraph 2
Replace-Text "This has made finding the high quality data more difficult and reduces the reusability of data in general" "This has made finding the high-quality data more difficult and has reduced the reusability of datasets in general"

# 3. "tests you might" -> "tests, you might"
Replace-Text "water quality tests you might" "water quality tests, you might"

# 4. "what the data are but there are" -> "what the data are about but there are"
Replace-Text "what the data are but there are still" "what the data are about but there are still"

# 5. "even missing" -> "despite missing"; "of high quality compared" -> "of very high quality compared"
Replace-Text "As it turns out, even missing those metadata elements, this dataset is of high quality compared to" "As it turns out, despite missing those metadata elements, this dataset is of very high quality compared to"

# 6. "ality metadata.  Data.wa.gov is not along; quality issued are" -> "...entries.  ...not alone; ..."
Replace-Text "quality metadata.  Data.wa.gov is not along; quality issued are" "quality metadata entries.  Data.wa.gov is not alone; quality issued are"

# 7. Insert "Cite that report" before the (non-highlighted surroundings) first "Kubler" citation, preserving its highlight run
$rng = $d.Content
$rng.Find.Execute("data portals (Kubler)") | Out-Null
$k = $d.Content
$k.Find.Execute("Kubler", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$k.InsertBefore("Cite that report")

# 10. "assessing the current state" -> "assess the current state"
Replace-Text "to understand data publishing behavior on the portal and assessing the current state of metadata quality on the portal." "to understand data publishing behavior on the portal and assess the current state of metadata quality on the portal."

# 13. "Any broad attempt to increase" -> "Any sweeping attempts to increase"
Replace-Text "Any broad attempt to increase metadata quality" "Any sweeping attempts to increase metadata quality"

# 14. "quality through completeness" -> "quality by examining completeness"
Replace-Text "I assessed metadata quality through completeness" "I assessed metadata quality by examining completeness"

Write-Output "done"
